$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$tr.Text = "https://replit.com/@HylandOutreach/RulesetExample"
$tr.Font.Size = 96
$tr.ActionSettings.Item(1).Hyperlink.Address = "https://replit.com/@HylandOutreach/RulesetExample"
